$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the title text in D12
$ws.Range("D12").Value = $null

# Update the link in E12 to the new URL
$ws.Range("E12").Value = "https://tensorflow.blog/2023/11/23/book-roadmap/"
